$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands and Aliasses")

# Apply the AutoFilter over the full table (header row 1 + data rows 2:56),
# filtering column A ("State") down to just "In Development" — this also
# recomputes which data rows are hidden by the filter.
$ws.Range("A1:E56").AutoFilter(1, @("In Development"), 7)

# Excel records the filter's range as a hidden, sheet-scoped defined name.
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "='Commands and Aliasses'!`$A`$1:`$E`$56")
$fd.Visible = $false

# Selecting the header row (as happens when a user clicks the column header
# after filtering) updates the saved selection/active-cell state.
$ws.Rows.Item(1).Select()
